$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.835.48"
$ws.Range("E2").Value = "  -6.24%  "
$ws.Range("D3").Value = "3.278.87"
$ws.Range("E3").Value = "  -8.35%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'177.27"
$ws.Range("E5").Value = "  -14.48%  "
$ws.Range("D6").Value = "'514.13"
$ws.Range("E6").Value = "  -8.94%  "
$ws.Range("D7").Value = "'0.589"
$ws.Range("E7").Value = "  -3.71%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.272.30"
$ws.Range("E8").Value = "  -8.35%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'0.611"
$ws.Range("E10").Value = "  -9.73%  "
$ws.Range("D11").Value = "'56.80"
$ws.Range("E11").Value = "  -6.48%  "
$ws.Range("D12").Value = "'0.130"
$ws.Range("E12").Value = "  -11.27%  "
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("E13").Value = "  -9.15%  "
$ws.Range("D14").Value = "'8.99"
$ws.Range("E14").Value = "  -12.03%  "
$ws.Range("D15").Value = "3.783.38"
$ws.Range("E15").Value = "  -8.83%  "
$ws.Range("E16").Value = "  -6.57%  "
$ws.Range("D17").Value = "3.258.34"
$ws.Range("E17").Value = "  -8.91%  "
$ws.Range("D18").Value = "63.437.30"
$ws.Range("E18").Value = "  -6.48%  "
$ws.Range("D19").Value = "'17.03"
$ws.Range("E19").Value = "  -9.66%  "
$ws.Range("D20").Value = "'10.73"
$ws.Range("E20").Value = "  -11.70%  "
$ws.Range("D21").Value = "'0.941"
$ws.Range("E21").Value = "  -11.06%  "
$ws.Range("D22").Value = "'367.33"
$ws.Range("E22").Value = "  -8.46%  "
$ws.Range("D23").Value = "'79.70"
$ws.Range("E23").Value = "  -5.32%  "
$ws.Range("D24").Value = "'3.61"
$ws.Range("E24").Value = "  -13.15%  "
$ws.Range("D25").Value = "'10.76"
$ws.Range("E25").Value = "  -13.82%  "
$ws.Range("D26").Value = "'3.79"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.62"
$ws.Range("E27").Value = "  -8.64%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'11.20"
$ws.Range("E28").Value = "  -9.63%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'8.23"
$ws.Range("E29").Value = "  -10.09%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'28.27"
$ws.Range("E30").Value = "  -9.87%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'633.56"
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.58"
$ws.Range("E32").Value = "  -14.37%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'10.99"
$ws.Range("E33").Value = "  -8.63%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'58.87"
$ws.Range("E34").Value = "  -6.79%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.102"
$ws.Range("E35").Value = "  -8.93%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'35.47"
$ws.Range("E37").Value = "  -13.34%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.371"
$ws.Range("E38").Value = "  -9.00%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.996"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.120"
$ws.Range("E40").Value = "  -9.11%  "
$ws.Range("D41").Value = "2.825.74"
$ws.Range("E41").Value = "  -10.73%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0647"
$ws.Range("E42").Value = "  -13.61%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  -19.23%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.57"
$ws.Range("E44").Value = "  -7.65%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  -14.62%  "
$ws.Range("D46").Value = "'0.0378"
$ws.Range("E46").Value = "  -7.65%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.123"
$ws.Range("E47").Value = "  -5.55%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.66"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'132.04"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'2.36"
$ws.Range("E51").Value = "  -8.95%  "
